{"js": "// Update the title paragraph (date) in place, preserving its run formatting.\nconst titleParagraph = context.document.body.paragraphs.getFirst();\ntitleParagraph.insertText(\"2025-08-23 Saturday\", Word.InsertLocation.replace);\n\n// Update every cell of the single answers table with the new equation text,\n// in one shot via Table.values -- this preserves each cell's existing\n// paragraph/run formatting (fonts, size, alignment) because only the <w:t>\n// text nodes are rewritten.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst newValues = [\n  [\"87-84=3\", \"35+10=45\", \"86+9=95\", \"1+30=31\", \"49+6=55\"],\n  [\"84-39=45\", \"22-6=16\", \"64+5=69\", \"62-39=23\", \"96-94=2\"],\n  [\"82-13=69\", \"85-73=12\", \"15+7=22\", \"72-18=54\", \"29-19=10\"],\n  [\"95-92=3\", \"74-16=58\", \"3+30=33\", \"31+24=55\", \"41+43=84\"],\n  [\"81-37=44\", \"27-18=9\", \"39-24=15\", \"19-17=2\", \"69+7=76\"],\n  [\"12+62=74\", \"81+2=83\", \"79+18=97\", \"56-25=31\", \"75-30=45\"],\n  [\"19+71=90\", \"93-86=7\", \"93-32=61\", \"63+4=67\", \"28-1=27\"],\n  [\"63-48=15\", \"30+11=41\", \"76-11=65\", \"91-61=30\", \"58-30=28\"],\n  [\"24+57=81\", \"41+58=99\", \"52-45=7\", \"87-40=47\", \"85-38=47\"],\n  [\"21+10=31\", \"11+49=60\", \"82-56=26\", \"75-51=24\", \"28+13=41\"],\n  [\"42+1=43\", \"73-72=1\", \"50-34=16\", \"7+68=75\", \"13+11=24\"],\n  [\"62-46=16\", \"33+62=95\", \"12+51=63\", \"41-11=30\", \"37+42=79\"],\n  [\"18+72=90\", \"22+39=61\", \"46+6=52\", \"20+79=99\", \"67+23=90\"],\n  [\"63-51=12\", \"50+44=94\", \"90-40=50\", \"56-54=2\", \"82+8=90\"],\n  [\"67-2=65\", \"19+45=64\", \"22+20=42\", \"25+25=50\", \"28+47=75\"],\n  [\"86-50=36\", \"28+32=60\", \"55+10=65\", \"50-46=4\", \"34-22=12\"],\n  [\"19-18=1\", \"75-13=62\", \"15+36=51\", \"42-17=25\", \"56-26=30\"],\n  [\"21-10=11\", \"6+66=72\", \"72-22=50\", \"36+50=86\", \"50+31=81\"],\n  [\"79-55=24\", \"52+36=88\", \"56-29=27\", \"21-13=8\", \"46-7=39\"],\n  [\"67-51=16\", \"62-58=4\", \"55-5=50\", \"78-76=2\", \"77+22=99\"],\n];\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the title paragraph (date) in place, preserving its run formatting.\n$d.Paragraphs.Item(1).Range.Text = '2025-08-23 Saturday'\n\n# Update every cell of the single answers table with the new equation text,\n# preserving each cell's existing paragraph/run formatting.\n$t = $d.Tables.Item(1)\n$newValues = @(\n    @('87-84=3', '35+10=45', '86+9=95', '1+30=31', '49+6=55'),\n    @('84-39=45', '22-6=16', '64+5=69', '62-39=23', '96-94=2'),\n    @('82-13=69', '85-73=12', '15+7=22', '72-18=54', '29-19=10'),\n    @('95-92=3', '74-16=58', '3+30=33', '31+24=55', '41+43=84'),\n    @('81-37=44', '27-18=9', '39-24=15', '19-17=2', '69+7=76'),\n    @('12+62=74', '81+2=83', '79+18=97', '56-25=31', '75-30=45'),\n    @('19+71=90', '93-86=7', '93-32=61', '63+4=67', '28-1=27'),\n    @('63-48=15', '30+11=41', '76-11=65', '91-61=30', '58-30=28'),\n    @('24+57=81', '41+58=99', '52-45=7', '87-40=47', '85-38=47'),\n    @('21+10=31', '11+49=60', '82-56=26', '75-51=24', '28+13=41'),\n    @('42+1=43', '73-72=1', '50-34=16', '7+68=75', '13+11=24'),\n    @('62-46=16', '33+62=95', '12+51=63', '41-11=30', '37+42=79'),\n    @('18+72=90', '22+39=61', '46+6=52', '20+79=99', '67+23=90'),\n    @('63-51=12', '50+44=94', '90-40=50', '56-54=2', '82+8=90'),\n    @('67-2=65', '19+45=64', '22+20=42', '25+25=50', '28+47=75'),\n    @('86-50=36', '28+32=60', '55+10=65', '50-46=4', '34-22=12'),\n    @('19-18=1', '75-13=62', '15+36=51', '42-17=25', '56-26=30'),\n    @('21-10=11', '6+66=72', '72-22=50', '36+50=86', '50+31=81'),\n    @('79-55=24', '52+36=88', '56-29=27', '21-13=8', '46-7=39'),\n    @('67-51=16', '62-58=4', '55-5=50', '78-76=2', '77+22=99')\n)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
